$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 19749.5
$ws.Range("J7").Value = 19749.5
$ws.Range("L7").Value = 19749.5
$ws.Range("N7").Value = -19973.5
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30586
$ws.Range("H14").Value = 19749.5
$ws.Range("J14").Value = 19749.5
$ws.Range("L14").Value = 19749.5
$ws.Range("N14").Value = -20131.5
$ws.Range("H19").Value = 2193928.2
$ws.Range("I19").Value = 5263834
$ws.Range("K19").Value = 5263834
$ws.Range("M19").Value = -5263659
$ws.Range("H33").Value = 250.92308
$ws.Range("I33").Value = 250.92308
$ws.Range("K33").Value = 250.92308
$ws.Range("M33").Value = -21.92308
$ws.Range("H43").Value = 2004.9
$ws.Range("I43").Value = 1315.25
$ws.Range("J43").Value = 2464.6667
$ws.Range("K43").Value = 1315.25
$ws.Range("L43").Value = 2464.6667
$ws.Range("M43").Value = -1246.25
$ws.Range("N43").Value = -2602.6667
$ws.Range("H129").Value = 850.47
$ws.Range("J129").Value = 870.3125
$ws.Range("L129").Value = 2610.9375
$ws.Range("N129").Value = -12610.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3868.0688
$ws.Range("I74").Value = 4169
$ws.Range("J74").Value = 3078.125
$ws.Range("K74").Value = 4169
$ws.Range("L74").Value = 3078.125
$ws.Range("M74").Value = -3295
$ws.Range("N74").Value = -4826.125
$ws.Range("H77").Value = 3868.0688
$ws.Range("I77").Value = 4169
$ws.Range("J77").Value = 3078.125
$ws.Range("K77").Value = 20845
$ws.Range("L77").Value = 15390.625
$ws.Range("M77").Value = -16477
$ws.Range("N77").Value = -24126.625
$ws.Range("H122").Value = 4577.385
$ws.Range("I122").Value = 1569.5
$ws.Range("J122").Value = 9390
$ws.Range("K122").Value = 4708.5
$ws.Range("L122").Value = 28170
$ws.Range("M122").Value = -2258.5
$ws.Range("N122").Value = -33070

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 3788.889
$ws.Range("I8").Value = 1314.2858
$ws.Range("J8").Value = 12450
$ws.Range("K8").Value = 1314.2858
$ws.Range("L8").Value = 12450
$ws.Range("M8").Value = -1174.2858
$ws.Range("N8").Value = -12730
$ws.Range("H86").Value = 2956.5715
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2956.5715
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H58").Value = 3101.1562
$ws.Range("I58").Value = 1797.5769
$ws.Range("J58").Value = 8750
$ws.Range("K58").Value = 1797.5769
$ws.Range("L58").Value = 8750
$ws.Range("M58").Value = -1594.5769
$ws.Range("N58").Value = -9156
$ws.Range("H68").Value = 99999
$ws.Range("J68").Value = 99999
$ws.Range("L68").Value = 99999
$ws.Range("N68").Value = -101497
$ws.Range("H71").Value = 99999
$ws.Range("J71").Value = 99999
$ws.Range("L71").Value = 299997
$ws.Range("N71").Value = -307485
$ws.Range("H136").Value = 3101.1562
$ws.Range("I136").Value = 1797.5769
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 5392.7307
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -2842.7307
$ws.Range("N136").Value = -31350
$ws.Range("H137").Value = 45324.285
$ws.Range("J137").Value = 45324.285
$ws.Range("L137").Value = 45324.285
$ws.Range("N137").Value = -55524.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1251
$ws.Range("I46").Value = 376.5
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 1129.5
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = -1038.5
$ws.Range("N46").Value = -9182
$ws.Range("H107").Value = 494.41666
$ws.Range("I107").Value = 354.3
$ws.Range("J107").Value = 1195
$ws.Range("K107").Value = 1062.9
$ws.Range("L107").Value = 3585
$ws.Range("M107").Value = 857.0999999999999
$ws.Range("N107").Value = -7425
$ws.Range("H113").Value = 3677069.8
$ws.Range("I113").Value = 612.82355
$ws.Range("J113").Value = 7353526.5
$ws.Range("K113").Value = 1838.47065
$ws.Range("L113").Value = 22060579.5
$ws.Range("M113").Value = 331.5293500000002
$ws.Range("N113").Value = -22064919.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 11251188
$ws.Range("I10").Value = 11251188
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 11251188
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -11251019
$ws.Range("N10").ClearContents()
$ws.Range("H120").Value = 26166.666
$ws.Range("J120").Value = 26166.666
$ws.Range("L120").Value = 26166.666
$ws.Range("N120").Value = -35842.666
$ws.Range("H122").Value = 9385.714
$ws.Range("J122").Value = 14250
$ws.Range("L122").Value = 42750
$ws.Range("N122").Value = -47650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3388.4092
$ws.Range("I7").Value = 2384.0625
$ws.Range("J7").Value = 6066.6665
$ws.Range("K7").Value = 2384.0625
$ws.Range("L7").Value = 6066.6665
$ws.Range("M7").Value = -2272.0625
$ws.Range("N7").Value = -6290.6665
$ws.Range("H22").Value = 2422.6365
$ws.Range("I22").Value = 1924.8334
$ws.Range("J22").Value = 3020
$ws.Range("K22").Value = 1924.8334
$ws.Range("L22").Value = 3020
$ws.Range("M22").Value = -1629.8334
$ws.Range("N22").Value = -3610
$ws.Range("H27").Value = 2422.6365
$ws.Range("I27").Value = 1924.8334
$ws.Range("J27").Value = 3020
$ws.Range("K27").Value = 1924.8334
$ws.Range("L27").Value = 3020
$ws.Range("M27").Value = -1817.8334
$ws.Range("N27").Value = -3234
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112
$ws.Range("H122").Value = 6502.5
$ws.Range("I122").Value = 2405
$ws.Range("J122").Value = 10600
$ws.Range("K122").Value = 7215
$ws.Range("L122").Value = 31800
$ws.Range("M122").Value = -4765
$ws.Range("N122").Value = -36700
$ws.Range("H126").Value = 3388.4092
$ws.Range("I126").Value = 2384.0625
$ws.Range("J126").Value = 6066.6665
$ws.Range("K126").Value = 7152.1875
$ws.Range("L126").Value = 18199.9995
$ws.Range("M126").Value = -4682.1875
$ws.Range("N126").Value = -23139.9995
$ws.Range("H132").Value = 6712.125
$ws.Range("I132").Value = 4100
$ws.Range("J132").Value = 9324.25
$ws.Range("K132").Value = 12300
$ws.Range("L132").Value = 27972.75
$ws.Range("M132").Value = -9770
$ws.Range("N132").Value = -33032.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4151.5884
$ws.Range("I122").Value = 1814.25
$ws.Range("J122").Value = 6229.222
$ws.Range("K122").Value = 5442.75
$ws.Range("L122").Value = 18687.666
$ws.Range("M122").Value = -2992.75
$ws.Range("N122").Value = -23587.666
$ws.Range("H126").Value = 628166.25
$ws.Range("I126").Value = 1436.6
$ws.Range("J126").Value = 1523494.2
$ws.Range("K126").Value = 4309.799999999999
$ws.Range("L126").Value = 4570482.6
$ws.Range("M126").Value = -1839.799999999999
$ws.Range("N126").Value = -4575422.6
$ws.Range("H132").Value = 13891203
$ws.Range("I132").Value = 1733.421
$ws.Range("J132").Value = 66671188
$ws.Range("K132").Value = 5200.263
$ws.Range("L132").Value = 200013564
$ws.Range("M132").Value = -2670.263
$ws.Range("N132").Value = -200018624
$ws.Range("H138").Value = 37975.6
$ws.Range("J138").Value = 37975.6
$ws.Range("L138").Value = 37975.6
$ws.Range("N138").Value = -48255.6
